$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, pushing existing rows 81..174 down to 82..175.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new record.
$ws.Range("A81").Value2 = 8
$ws.Range("B81").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C81").Value2 = "Coquimbo"
$ws.Range("D81").Value2 = 44650
$ws.Range("D81").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E81").Value2 = 4
$ws.Range("F81").Value2 = 100112037
$ws.Range("G81").Value2 = "Cebollín"
$ws.Range("H81").Value2 = "Sin especificar"
$ws.Range("I81").Value2 = "Primera"
$ws.Range("J81").Value2 = 1800
$ws.Range("K81").Value2 = 5500
$ws.Range("L81").Value2 = 6000
$ws.Range("M81").Value2 = 5750
$ws.Range("N81").Value2 = "$/paquete 36 unidades"
$ws.Range("O81").Value2 = "Provincia del Elquí"
$ws.Range("P81").Value2 = 160
$ws.Range("Q81").Value2 = 36
$ws.Range("R81").Value2 = "Hortaliza"
